$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 12, shifting the existing
# rows 12..82 down to 14..84 (formatting/styles are carried along).
$ws.Rows("12:13").Insert()

# Populate the two newly inserted rows with the new weekly records.
# Row 12: Lapins / Primera
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(12, 3).Value = "Bíobío"
$ws.Cells.Item(12, 4).Value = 44558
$ws.Cells.Item(12, 5).Value = 8
$ws.Cells.Item(12, 6).Value = "Fruta"
$ws.Cells.Item(12, 7).Value = 100103
$ws.Cells.Item(12, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(12, 9).Value = 100103001
$ws.Cells.Item(12, 10).Value = "Cereza"
$ws.Cells.Item(12, 11).Value = "Lapins"
$ws.Cells.Item(12, 12).Value = "Primera"
$ws.Cells.Item(12, 13).Value = 100
$ws.Cells.Item(12, 14).Value = 4500
$ws.Cells.Item(12, 15).Value = 5000
$ws.Cells.Item(12, 16).Value = 4750
$ws.Cells.Item(12, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(12, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(12, 19).Value = 475
$ws.Cells.Item(12, 20).Value = 10

# Row 13: Santina / Primera
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(13, 3).Value = "Bíobío"
$ws.Cells.Item(13, 4).Value = 44558
$ws.Cells.Item(13, 5).Value = 8
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100103
$ws.Cells.Item(13, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(13, 9).Value = 100103001
$ws.Cells.Item(13, 10).Value = "Cereza"
$ws.Cells.Item(13, 11).Value = "Santina"
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 200
$ws.Cells.Item(13, 14).Value = 4500
$ws.Cells.Item(13, 15).Value = 5000
$ws.Cells.Item(13, 16).Value = 4750
$ws.Cells.Item(13, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(13, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(13, 19).Value = 475
$ws.Cells.Item(13, 20).Value = 10
